$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 38. This shifts existing rows 38-90 down to 39-91,
# carrying all of their D/J/K/L/M/P (and other) values with them intact.
$ws.Rows("38").Insert()

# Populate the newly inserted row 38 with a new weekly record. The
# "template" columns (A,B,C,E,F,G,H,I,N,O,Q,R) match the surrounding rows
# for this Terminal/Cilantro series; K,L,M,P repeat the same price-range
# figures as before, while D (date) and J (price) carry the new values.
$ws.Range("A38").Value = 8
$ws.Range("B38").Value = "Terminal La Palmera de La Serena"
$ws.Range("C38").Value = "Coquimbo"
$ws.Range("D38").Value2 = 44495
$ws.Range("E38").Value = 4
$ws.Range("F38").Value = 100112040
$ws.Range("G38").Value = "Cilantro"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 3200
$ws.Range("K38").Value = 1300
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = 1400
$ws.Range("N38").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O38").Value = "Provincia del Elquí"
$ws.Range("P38").Value = 933
$ws.Range("Q38").Value = 1.5
$ws.Range("R38").Value = "Hortaliza"
